# Sync Local VIP Card BOM
#
# The sheet that held the AVI-ELF II VIP card BOM is renamed to reflect the
# "Digikey.ca" cart it was generated from, and the saved cursor/selection
# moves from the old working cell (F27) back to the top of the table (A3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Make sure we're working on the (only) sheet before touching its view state.
$ws.Activate()

# Rename the worksheet tab: "AVIELF2-VIP" -> "Digikey.ca"
$ws.Name = "Digikey.ca"

# Move the saved selection/active cell from F27 to A3.
$ws.Range("A3").Select()
